# Logged Week 15 and simulated Week 16
# Update the "H" row (row 2) passing-depth stats on both the OFF and DEF
# sheets to reflect the newly logged/simulated week.

$wb = $excel.ActiveWorkbook

# --- OFF sheet (Short Att, Short Comp, Deep Att, Deep Comp, Short Int) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 415
$wsOff.Range("C2").Value = 272
$wsOff.Range("D2").Value = 94
$wsOff.Range("E2").Value = 29
$wsOff.Range("F2").Value = 4

# --- DEF sheet (Short Att, Short Comp, Deep Att, Deep Comp, Short Int) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 523
$wsDef.Range("C2").Value = 371
$wsDef.Range("D2").Value = 114
$wsDef.Range("E2").Value = 56
$wsDef.Range("F2").Value = 11
